$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new trade row (row 10), copying formats from row 9 above it so that
# styled cells (date column A and boolean column G) keep the same style index.
$r = 10
$prev = $r - 1

$ws.Cells.Item($prev, 1).Copy($ws.Cells.Item($r, 1))
$ws.Cells.Item($prev, 7).Copy($ws.Cells.Item($r, 7))

$ws.Cells.Item($r, 1).Value2 = 42654.743831018517
$ws.Cells.Item($r, 2).Value2 = $true
$ws.Cells.Item($r, 3).Value2 = 10185.17
$ws.Cells.Item($r, 4).Value2 = 10012.950000000001
$ws.Cells.Item($r, 5).Value2 = 18.870000999999998
$ws.Cells.Item($r, 6).Value2 = 19.52
$ws.Cells.Item($r, 7).Value2 = $false
$ws.Cells.Item($r, 8).Value2 = 3.44
$ws.Cells.Item($r, 9).Value2 = $false

$ws.Columns.Item(1).ColumnWidth = 14.541666666666666
$ws.Columns.Item(2).ColumnWidth = 7.416666666666667
$ws.Columns.Item(3).ColumnWidth = 8.041666666666666
$ws.Columns.Item(4).ColumnWidth = 10.416666666666666
$ws.Columns.Item(5).ColumnWidth = 9.041666666666666
$ws.Columns.Item(6).ColumnWidth = 6.166666666666667
$ws.Columns.Item(7).ColumnWidth = 9.541666666666666
$ws.Columns.Item(8).ColumnWidth = 13.791666666666666
$ws.Columns.Item(9).ColumnWidth = 11.041666666666666
